# Add a new "namInUrdu" column after "name" (inserted as column C),
# shifting productImage/batchCode/... one column to the right, and
# populate the new column with "urdu" for the existing product rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (shifts existing C:I to D:J)
$ws.Columns.Item(3).Insert()

# Header for the new column
$ws.Range("C1").Value = "namInUrdu"

# Data for the new column
$ws.Range("C2").Value = "urdu"
$ws.Range("C3").Value = "urdu"

# Match the resulting selection/view state
$ws.Range("C5").Select()
